# Applies the weekly CompStat data refresh described in the commit:
# "New crime data collected" -- updates header Volume/Number + date range,
# and refreshes the crime-complaints table (rows 15-31) with the new weeks figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and date range) ---
$ws.Range("A8").Value = "Volume 32   Number  29"
$ws.Range("C9").Value = "Report Covering the Week  7/14/2025  Through  7/20/2025"

# --- Data table updates (rows 15-31) ---
# Row 15
$ws.Range("C15").Value = 1
$ws.Range("C15").NumberFormat = "#,##0"
$ws.Range("F15").Value = 1
$ws.Range("F15").NumberFormat = "#,##0"
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 9
$ws.Range("K15").Value = 80
$ws.Range("L15").Value = 350
$ws.Range("M15").Value = 800
$ws.Range("N15").Value = 0
# Row 16
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = 100
$ws.Range("F16").Value = 25
$ws.Range("G16").Value = 13
$ws.Range("H16").Value = 92.307692307692
$ws.Range("I16").Value = 129
$ws.Range("J16").Value = 133
$ws.Range("K16").Value = -3.007518796992
$ws.Range("L16").Value = 69.736842105263
$ws.Range("M16").Value = 17.272727272727
$ws.Range("N16").Value = -80.717488789237
# Row 17
$ws.Range("C17").Value = 2
$ws.Range("D17").Value = 10
$ws.Range("E17").Value = -80
$ws.Range("F17").Value = 27
$ws.Range("G17").Value = 35
$ws.Range("H17").Value = -22.857142857142
$ws.Range("I17").Value = 194
$ws.Range("J17").Value = 192
$ws.Range("K17").Value = 1.041666666666
$ws.Range("L17").Value = 36.619718309859
$ws.Range("M17").Value = 212.903225806452
$ws.Range("N17").Value = -11.818181818181
# Row 18
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 6
$ws.Range("E18").Value = -83.333333333333
$ws.Range("F18").Value = 12
$ws.Range("G18").Value = 18
$ws.Range("H18").Value = -33.333333333333
$ws.Range("I18").Value = 114
$ws.Range("J18").Value = 119
$ws.Range("K18").Value = -4.201680672268
$ws.Range("L18").Value = -10.236220472440
$ws.Range("M18").Value = 78.125
$ws.Range("N18").Value = -74.439461883408
# Row 19
$ws.Range("C19").Value = 14
$ws.Range("D19").Value = 13
$ws.Range("E19").Value = 7.692307692307
$ws.Range("F19").Value = 48
$ws.Range("G19").Value = 59
$ws.Range("H19").Value = -18.644067796610
$ws.Range("I19").Value = 369
$ws.Range("J19").Value = 407
$ws.Range("K19").Value = -9.336609336609
$ws.Range("L19").Value = -3.90625
$ws.Range("M19").Value = 45.849802371541
$ws.Range("N19").Value = -35.714285714285
# Row 20
$ws.Range("C20").Value = 3
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = -25
$ws.Range("F20").Value = 13
$ws.Range("G20").Value = 14
$ws.Range("H20").Value = -7.142857142857
$ws.Range("I20").Value = 56
$ws.Range("J20").Value = 46
$ws.Range("K20").Value = 21.739130434782
$ws.Range("L20").Value = 40
$ws.Range("M20").Value = 107.407407407407
$ws.Range("N20").Value = -85.822784810126
# Row 21
$ws.Range("C21").Value = 27
$ws.Range("D21").Value = 36
$ws.Range("E21").Value = -25
$ws.Range("F21").Value = 126
$ws.Range("G21").Value = 141
$ws.Range("H21").Value = -10.638297872340
$ws.Range("I21").Value = 871
$ws.Range("J21").Value = 904
$ws.Range("K21").Value = -3.650442477876
$ws.Range("L21").Value = 12.823834196891
$ws.Range("M21").Value = 68.146718146718
$ws.Range("N21").Value = -62.408286577470
# Row 22
$ws.Range("C22").Value = 1
$ws.Range("D22").Value = 1
$ws.Range("E22").Value = 0
$ws.Range("I22").Value = 35
$ws.Range("J22").Value = 27
$ws.Range("K22").Value = 29.629629629629
$ws.Range("L22").Value = 2.941176470588
$ws.Range("M22").Value = -7.894736842105
# Row 23
$ws.Range("D23").Value = 1
$ws.Range("D23").NumberFormat = "#,##0"
$ws.Range("E23").Value = -100
$ws.Range("E23").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("F23").Value = 3
$ws.Range("G23").Value = 2
$ws.Range("H23").Value = 50
$ws.Range("J23").Value = 24
$ws.Range("K23").Value = -16.666666666666
$ws.Range("M23").Value = 53.846153846153
# Row 24
$ws.Range("C24").Value = 38
$ws.Range("D24").Value = 48
$ws.Range("E24").Value = -20.833333333333
$ws.Range("F24").Value = 148
$ws.Range("G24").Value = 187
$ws.Range("H24").Value = -20.855614973262
$ws.Range("I24").Value = 918
$ws.Range("J24").Value = 1202
$ws.Range("K24").Value = -23.627287853577
$ws.Range("L24").Value = -17.889087656529
$ws.Range("M24").Value = 12.915129151291
# Row 25
$ws.Range("C25").Value = 32
$ws.Range("D25").Value = 42
$ws.Range("E25").Value = -23.809523809523
$ws.Range("F25").Value = 126
$ws.Range("G25").Value = 171
$ws.Range("H25").Value = -26.315789473684
$ws.Range("I25").Value = 713
$ws.Range("J25").Value = 1089
$ws.Range("K25").Value = -34.527089072543
$ws.Range("L25").Value = -29.405940594059
# Row 26
$ws.Range("C26").Value = 11
$ws.Range("D26").Value = 4
$ws.Range("E26").Value = 175
$ws.Range("F26").Value = 46
$ws.Range("G26").Value = 29
$ws.Range("H26").Value = 58.620689655172
$ws.Range("I26").Value = 240
$ws.Range("J26").Value = 259
$ws.Range("K26").Value = -7.335907335907
$ws.Range("L26").Value = 17.073170731707
$ws.Range("M26").Value = 22.448979591836
# Row 27
$ws.Range("C27").Value = 1
$ws.Range("C27").NumberFormat = "#,##0"
$ws.Range("F27").Value = 1
$ws.Range("F27").NumberFormat = "#,##0"
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 13
$ws.Range("K27").Value = 44.444444444444
$ws.Range("L27").Value = 160
# Row 28
$ws.Range("C28").Value = 1
$ws.Range("E28").Value = -50
$ws.Range("F28").Value = 8
$ws.Range("G28").Value = 7
$ws.Range("H28").Value = 14.285714285714
$ws.Range("I28").Value = 50
$ws.Range("J28").Value = 31
$ws.Range("K28").Value = 61.290322580645
$ws.Range("L28").Value = 72.413793103448
# Row 31
$ws.Range("D31").Value = 2
$ws.Range("D31").NumberFormat = "#,##0"
$ws.Range("E31").Value = -100
$ws.Range("E31").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("G31").Value = 2
$ws.Range("G31").NumberFormat = "#,##0"
$ws.Range("H31").Value = -100
$ws.Range("H31").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("J31").Value = 8
$ws.Range("K31").Value = -25
